$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Replace-InCell($table, $row, $col, $old, $new) {
    $cellRng = $table.Cell($row, $col).Range
    $rng = $d.Range($cellRng.Start, $cellRng.End)
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 0, $false, $new, 1) | Out-Null
}

# CNPJ number
Replace-InCell $t 2 1 "1106462000199" "312"

# I - Revenda de mercadorias com dispensa de emissao de documento fiscal
Replace-InCell $t 6 2 "1000.00" "42.00"

# II - Revenda de mercadorias com documento fiscal emitido
Replace-InCell $t 7 2 "0.00" "42.00"

# III - Total das receitas com revenda de mercadorias (I + II)
Replace-InCell $t 8 2 "1000.00" "84.00"

# IV - Venda de produtos industrializados com dispensa de emissao de documento fiscal
Replace-InCell $t 10 2 "0.00" "41.00"

# V - Venda de produtos industrializados com documento fiscal emitido
Replace-InCell $t 11 2 "0.00" "41.00"

# VI - Total das receitas com venda de produtos industrializados (IV + V)
Replace-InCell $t 12 2 "0.00" "82.00"

# VII - Receita com prestacao de servicos com dispensa de emissao de documento fiscal
Replace-InCell $t 14 2 "0.00" "1.00"

# VIII - Receita com prestacao de servicos com documento fiscal emitido
Replace-InCell $t 15 2 "0.00" "1.00"

# IX - Total das receitas com prestacao de servicos (VII + VIII)
Replace-InCell $t 16 2 "0.00" "2.00"

# X - Total geral das receitas brutas no mes (III + VI + IX)
Replace-InCell $t 17 2 "1000.00" "168.00"
